$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.948.35'
$ws.Range("D3").Value = '1.673.03'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("D5").Value = '214.86'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").Value = '20.18'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").Value = '1.907.92'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Value = '1.706.06'
$ws.Range("E13").Value = '  +3.08%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '65.66'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '26.946.66'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '8.09'
$ws.Range("E18").Value = '  +3.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '234.50'
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").Value = '0.0₃0732'
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = '4.44'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '9.16'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '145.65'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '7.16'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").Value = '15.98'
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("D33").Value = '1.474.99'
$ws.Range("E33").Value = '  -5.13%  '
$ws.Range("E34").Value = '  +2.25%  '
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("D36").Value = '2.41'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").Value = '0.578'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").Value = '0.896'
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("E40").Value = '  +7.66%  '
$ws.Range("D41").Value = '5.81'
$ws.Range("E41").Value = '  -3.77%  '
$ws.Range("D43").Value = '2.29'
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("D44").Value = '66.73'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D45").Value = '1.814.73'
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = '90.41'
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '7.69'
$ws.Range("E51").Value = '  +0.56%  '

Write-Output "Applied 76 cell updates"
